# Updates cryptos.xlsx price (D) and 1h-volume-change (E) columns with the
# latest scraped values from the GitHub Actions run.
# Columns D/E hold plain text (e.g. "26.153.42", "  -4.32%  ") rather than
# numbers, so we force a Text number format before writing any value that
# could otherwise be auto-parsed by Excel as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.153.42"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.08"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.00"
$ws.Range("E5").Value = "  -3.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5101"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2588"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06405"
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.97"
$ws.Range("E10").Value = "  -4.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07815"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.657.48"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.291"
$ws.Range("E13").Value = "  -4.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.884.65"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5536"
$ws.Range("E15").Value = "  -4.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8021"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.14"
$ws.Range("E17").Value = "  -5.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.192.35"
$ws.Range("E18").Value = "  -4.19%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "209.49"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  -3.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.024"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.77"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.728"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1169"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.979"
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.78"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05102"
$ws.Range("E30").Value = "  -4.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.246"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.351"
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.228"
$ws.Range("E33").Value = "  -5.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.562"
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9284"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5708"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.153.63"
$ws.Range("E39").Value = "  +8.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01591"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8371"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.650"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.53"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.794.78"
$ws.Range("E45").Value = "  -3.10%  "
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4549"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.81"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.884"
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("E51").Value = "  -3.38%  "
